$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)

# Resize / reposition the content placeholder (values given in points;
# EMU targets: off x=838200 y=1333568, ext cx=10515600 cy=4033562)
$shp.Left   = 66.0
$shp.Top    = 105.00535433070866
$shp.Width  = 828.0
$shp.Height = 317.603318

$tr = $shp.TextFrame.TextRange

$para1 = "La communication asynchrone existe depuis longtemps et sont adopt" + [char]0x00E9 + "s dans plusieurs architectures de communication surtout dans les IoT. Elle est diff" + [char]0x00E9 + "rentes de celle de synchrone aux niveaux du couplage, de protocole, de l" + [char]0x2019 + [char]0x00E9 + "volutivit" + [char]0x00E9 + ", l" + [char]0x2019 + "absence de standard, des MOM (Message-Oriented Middleware) et surtout de gestion d" + [char]0x2019 + "erreur quant on parle des APIs."

$para2 = "Ainsi, nous d" + [char]0x00E9 + "vons retenir que les APIs HTTP ou APIs pub/sub n" + [char]0x2019 + "est qu" + [char]0x2019 + "une question de communication synchrone ou asynchrone. Mais le plus important est que nous ne devons jamais oublier les consommateurs des ces APIs. Ce qui implique une sp" + [char]0x00E9 + "cification l" + [char]0x00E9 + "ger, adapt" + [char]0x00E9 + " et qui peut " + [char]0x00EA + "tre am" + [char]0x00E9 + "lior" + [char]0x00E9 + " au fur et " + [char]0x00E0 + " mesure l" + [char]0x2019 + "avancer technologique et d" + [char]0x2019 + "ajout des contrainte."

$tr.Text = $para1 + [char]13 + $para2

# Bold run for the "MOM" acronym in paragraph 1.
$momRange = $tr.Characters(264, 3)
$momRange.Font.Bold = $true

# Force a dedicated run for "Oriented" (mirrors the author's spell-checker
# generated run split) without altering its visible formatting.
$orientedRange = $tr.Characters(277, 8)
$orientedRange.Font.Size = $orientedRange.Font.Size
